# ---------------------------------------------------------------------------
# "Refined metadata to be additional tab"
#
# 1. Add a new "metadata" worksheet positioned after the existing "data"
#    sheet, re-activate "data" so the active tab is unchanged.
# 2. Populate the "metadata" sheet with a header row + one data row
#    describing the PanelApp query used to build the "data" sheet,
#    re-using the same bold/bordered header style as the "data" sheet.
# 3. Refresh the "time_taken" timestamps (column F) on the "data" sheet to
#    reflect the re-run query time.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- 1. Insert the "metadata" worksheet right after "data" ----------------
$meta = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$meta.Name = "metadata"
$ws.Activate()

# --- 2. Populate "metadata" ------------------------------------------------

# Header row (B1:G1) - reuse the bold/bordered/centered style already used
# for the "data" sheet's own header row.
$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)

# Index cell (A2) - reuse the same style as data!A2.
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Retinitis Pigmentosa Superpanel"
$meta.Range("C2").Value = 3100

# "data_version" must stay a text value ("0.132"), not be coerced to a
# number - enter it as a formula producing the text string, then
# collapse the formula down to its literal value so the cell ends up a
# plain, unstyled text cell (no stray "text number format" left behind).
$meta.Range("D2").Formula = "=""0.132"""
$meta.Range("D2").Copy()
$meta.Range("D2").PasteSpecial(-4163)

$meta.Range("E2").Value = "2021-09-18T08:15:22.721768Z"
$meta.Range("F2").Value = "2021-10-05 14:35:39.838765"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3100/?format=json"

# --- 3. Refresh "data" sheet timestamps (column F) --------------------------
$ws.Range("F2").Value = "2021-10-05 14:35:39.842251"
$ws.Range("F3").Value = "2021-10-05 14:35:39.842258"
$ws.Range("F4").Value = "2021-10-05 14:35:39.842262"
$ws.Range("F5").Value = "2021-10-05 14:35:39.842264"
$ws.Range("F6").Value = "2021-10-05 14:35:39.842267"
$ws.Range("F7").Value = "2021-10-05 14:35:39.842270"
$ws.Range("F8").Value = "2021-10-05 14:35:39.842272"
$ws.Range("F9").Value = "2021-10-05 14:35:39.842275"
$ws.Range("F10").Value = "2021-10-05 14:35:39.842278"
$ws.Range("F11").Value = "2021-10-05 14:35:39.842281"
$ws.Range("F12").Value = "2021-10-05 14:35:39.842283"
$ws.Range("F13").Value = "2021-10-05 14:35:39.842286"
$ws.Range("F14").Value = "2021-10-05 14:35:39.842289"
$ws.Range("F15").Value = "2021-10-05 14:35:39.842291"
$ws.Range("F16").Value = "2021-10-05 14:35:39.842294"
$ws.Range("F17").Value = "2021-10-05 14:35:39.842296"
$ws.Range("F18").Value = "2021-10-05 14:35:39.842299"
$ws.Range("F19").Value = "2021-10-05 14:35:39.842302"
$ws.Range("F20").Value = "2021-10-05 14:35:39.842305"
$ws.Range("F21").Value = "2021-10-05 14:35:39.842308"
$ws.Range("F22").Value = "2021-10-05 14:35:39.842310"
$ws.Range("F23").Value = "2021-10-05 14:35:39.842313"
$ws.Range("F24").Value = "2021-10-05 14:35:39.842316"
$ws.Range("F25").Value = "2021-10-05 14:35:39.842318"
$ws.Range("F26").Value = "2021-10-05 14:35:39.842321"
$ws.Range("F27").Value = "2021-10-05 14:35:39.842324"
$ws.Range("F28").Value = "2021-10-05 14:35:39.842327"
$ws.Range("F29").Value = "2021-10-05 14:35:39.842329"
$ws.Range("F30").Value = "2021-10-05 14:35:39.842332"
$ws.Range("F31").Value = "2021-10-05 14:35:39.842334"
$ws.Range("F32").Value = "2021-10-05 14:35:39.842337"
$ws.Range("F33").Value = "2021-10-05 14:35:39.842340"
$ws.Range("F34").Value = "2021-10-05 14:35:39.842343"
$ws.Range("F35").Value = "2021-10-05 14:35:39.842345"
$ws.Range("F36").Value = "2021-10-05 14:35:39.842348"
$ws.Range("F37").Value = "2021-10-05 14:35:39.842351"
$ws.Range("F38").Value = "2021-10-05 14:35:39.842353"
$ws.Range("F39").Value = "2021-10-05 14:35:39.842356"
$ws.Range("F40").Value = "2021-10-05 14:35:39.842359"
$ws.Range("F41").Value = "2021-10-05 14:35:39.842361"
$ws.Range("F42").Value = "2021-10-05 14:35:39.842365"
$ws.Range("F43").Value = "2021-10-05 14:35:39.842367"
$ws.Range("F44").Value = "2021-10-05 14:35:39.842370"
$ws.Range("F45").Value = "2021-10-05 14:35:39.842373"
$ws.Range("F46").Value = "2021-10-05 14:35:39.842375"
$ws.Range("F47").Value = "2021-10-05 14:35:39.842378"
$ws.Range("F48").Value = "2021-10-05 14:35:39.842381"
$ws.Range("F49").Value = "2021-10-05 14:35:39.842383"
$ws.Range("F50").Value = "2021-10-05 14:35:39.842386"
$ws.Range("F51").Value = "2021-10-05 14:35:39.842389"
$ws.Range("F52").Value = "2021-10-05 14:35:39.842391"
$ws.Range("F53").Value = "2021-10-05 14:35:39.842394"
$ws.Range("F54").Value = "2021-10-05 14:35:39.842397"
$ws.Range("F55").Value = "2021-10-05 14:35:39.842400"
$ws.Range("F56").Value = "2021-10-05 14:35:39.842402"
$ws.Range("F57").Value = "2021-10-05 14:35:39.842405"
$ws.Range("F58").Value = "2021-10-05 14:35:39.842407"
$ws.Range("F59").Value = "2021-10-05 14:35:39.842410"
$ws.Range("F60").Value = "2021-10-05 14:35:39.842413"
$ws.Range("F61").Value = "2021-10-05 14:35:39.842416"
$ws.Range("F62").Value = "2021-10-05 14:35:39.842418"
$ws.Range("F63").Value = "2021-10-05 14:35:39.842421"
$ws.Range("F64").Value = "2021-10-05 14:35:39.842424"
$ws.Range("F65").Value = "2021-10-05 14:35:39.842426"
$ws.Range("F66").Value = "2021-10-05 14:35:39.842431"
$ws.Range("F67").Value = "2021-10-05 14:35:39.842434"
$ws.Range("F68").Value = "2021-10-05 14:35:39.842437"
$ws.Range("F69").Value = "2021-10-05 14:35:39.842439"
$ws.Range("F70").Value = "2021-10-05 14:35:39.842442"
$ws.Range("F71").Value = "2021-10-05 14:35:39.842444"
$ws.Range("F72").Value = "2021-10-05 14:35:39.842447"
$ws.Range("F73").Value = "2021-10-05 14:35:39.842450"
$ws.Range("F74").Value = "2021-10-05 14:35:39.842452"
$ws.Range("F75").Value = "2021-10-05 14:35:39.842455"
$ws.Range("F76").Value = "2021-10-05 14:35:39.842457"
$ws.Range("F77").Value = "2021-10-05 14:35:39.842460"
$ws.Range("F78").Value = "2021-10-05 14:35:39.842465"
$ws.Range("F79").Value = "2021-10-05 14:35:39.842468"
$ws.Range("F80").Value = "2021-10-05 14:35:39.842471"
$ws.Range("F81").Value = "2021-10-05 14:35:39.842473"
$ws.Range("F82").Value = "2021-10-05 14:35:39.842476"
$ws.Range("F83").Value = "2021-10-05 14:35:39.842479"
$ws.Range("F84").Value = "2021-10-05 14:35:39.842481"
$ws.Range("F85").Value = "2021-10-05 14:35:39.842484"
$ws.Range("F86").Value = "2021-10-05 14:35:39.842486"
$ws.Range("F87").Value = "2021-10-05 14:35:39.842489"
$ws.Range("F88").Value = "2021-10-05 14:35:39.842492"
$ws.Range("F89").Value = "2021-10-05 14:35:39.842494"
$ws.Range("F90").Value = "2021-10-05 14:35:39.842497"
$ws.Range("F91").Value = "2021-10-05 14:35:39.842500"
$ws.Range("F92").Value = "2021-10-05 14:35:39.842502"
$ws.Range("F93").Value = "2021-10-05 14:35:39.842505"
$ws.Range("F94").Value = "2021-10-05 14:35:39.842509"
$ws.Range("F95").Value = "2021-10-05 14:35:39.842512"
$ws.Range("F96").Value = "2021-10-05 14:35:39.842514"
$ws.Range("F97").Value = "2021-10-05 14:35:39.842517"
$ws.Range("F98").Value = "2021-10-05 14:35:39.842519"
$ws.Range("F99").Value = "2021-10-05 14:35:39.842522"
$ws.Range("F100").Value = "2021-10-05 14:35:39.842524"
$ws.Range("F101").Value = "2021-10-05 14:35:39.842527"
$ws.Range("F102").Value = "2021-10-05 14:35:39.842529"
$ws.Range("F103").Value = "2021-10-05 14:35:39.842532"
$ws.Range("F104").Value = "2021-10-05 14:35:39.842534"
$ws.Range("F105").Value = "2021-10-05 14:35:39.842537"
$ws.Range("F106").Value = "2021-10-05 14:35:39.842539"
$ws.Range("F107").Value = "2021-10-05 14:35:39.842542"
$ws.Range("F108").Value = "2021-10-05 14:35:39.842544"
$ws.Range("F109").Value = "2021-10-05 14:35:39.842547"
$ws.Range("F110").Value = "2021-10-05 14:35:39.842551"
$ws.Range("F111").Value = "2021-10-05 14:35:39.842554"
$ws.Range("F112").Value = "2021-10-05 14:35:39.842557"
$ws.Range("F113").Value = "2021-10-05 14:35:39.842559"
$ws.Range("F114").Value = "2021-10-05 14:35:39.842562"
$ws.Range("F115").Value = "2021-10-05 14:35:39.842564"
$ws.Range("F116").Value = "2021-10-05 14:35:39.842567"
$ws.Range("F117").Value = "2021-10-05 14:35:39.842569"
$ws.Range("F118").Value = "2021-10-05 14:35:39.842572"
$ws.Range("F119").Value = "2021-10-05 14:35:39.842575"
$ws.Range("F120").Value = "2021-10-05 14:35:39.842577"
$ws.Range("F121").Value = "2021-10-05 14:35:39.842580"
$ws.Range("F122").Value = "2021-10-05 14:35:39.842582"
$ws.Range("F123").Value = "2021-10-05 14:35:39.842585"
$ws.Range("F124").Value = "2021-10-05 14:35:39.842588"
$ws.Range("F125").Value = "2021-10-05 14:35:39.842590"
$ws.Range("F126").Value = "2021-10-05 14:35:39.842593"
$ws.Range("F127").Value = "2021-10-05 14:35:39.842595"
$ws.Range("F128").Value = "2021-10-05 14:35:39.842598"
$ws.Range("F129").Value = "2021-10-05 14:35:39.842600"
$ws.Range("F130").Value = "2021-10-05 14:35:39.842605"
$ws.Range("F131").Value = "2021-10-05 14:35:39.842607"
$ws.Range("F132").Value = "2021-10-05 14:35:39.842611"
$ws.Range("F133").Value = "2021-10-05 14:35:39.842613"
$ws.Range("F134").Value = "2021-10-05 14:35:39.842616"
$ws.Range("F135").Value = "2021-10-05 14:35:39.842618"
$ws.Range("F136").Value = "2021-10-05 14:35:39.842621"
$ws.Range("F137").Value = "2021-10-05 14:35:39.842624"
$ws.Range("F138").Value = "2021-10-05 14:35:39.842626"
$ws.Range("F139").Value = "2021-10-05 14:35:39.842629"
$ws.Range("F140").Value = "2021-10-05 14:35:39.842632"

Write-Host "metadata tab added; timestamps refreshed"
